$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"23.2600000000002"
$ws.Range("H2").Value = [double]"7.123612810744362e-10"
$ws.Range("I2").Value = [double]"7.123612810744362e-10"
$ws.Range("L2").Value = [double]"58.29426837345365"
$ws.Range("M2").Value = "[43.0894447422762, 73.49909200463111]"
$ws.Range("N2").Value = [double]"8.779499349742537e-10"
$ws.Range("O2").Value = [double]"8.779499349742537e-10"
$ws.Range("P2").Value = [double]"1.591237119836271"
$ws.Range("Q2").Value = "[1.2893423303021168, 1.8931319093704255]"
$ws.Range("R2").Value = [double]"7.72715225139109e-14"
$ws.Range("S2").Value = [double]"7.72715225139109e-14"
$ws.Range("T2").Value = [double]"56.66828097303738"
$ws.Range("U2").Value = "[46.87793292165064, 66.45862902442413]"
$ws.Range("V2").Value = [double]"3.552713678800501e-15"
$ws.Range("W2").Value = [double]"3.552713678800501e-15"
$ws.Range("X2").Value = [double]"17.36932932932948"
$ws.Range("Y2").Value = [double]"16.25173173173187"
$ws.Range("Z2").Value = [double]"18.48692692692709"
$ws.Range("F3").Value = [double]"23.2600000000002"
$ws.Range("H3").Value = [double]"5.005255432344313e-09"
$ws.Range("I3").Value = [double]"5.005255432344313e-09"
$ws.Range("L3").Value = [double]"56.0359478211881"
$ws.Range("M3").Value = "[38.73432814097556, 73.33756750140064]"
$ws.Range("N3").Value = [double]"5.206607367469474e-08"
$ws.Range("O3").Value = [double]"5.206607367469474e-08"
$ws.Range("P3").Value = [double]"1.817658211986887"
$ws.Range("Q3").Value = "[1.4654476241970409, 2.169868799776734]"
$ws.Range("R3").Value = [double]"1.527666881884215e-13"
$ws.Range("S3").Value = [double]"1.527666881884215e-13"
$ws.Range("T3").Value = [double]"54.75472281454606"
$ws.Range("U3").Value = "[44.61281646811304, 64.89662916097909]"
$ws.Range("V3").Value = [double]"3.530509218307998e-14"
$ws.Range("W3").Value = [double]"3.530509218307998e-14"
$ws.Range("X3").Value = [double]"16.53113113113127"
$ws.Range("Y3").Value = [double]"15.2272672672674"
$ws.Range("Z3").Value = [double]"17.83499499499515"
$ws.Range("F4").Value = [double]"23.2600000000002"
$ws.Range("H4").Value = [double]"2.754685368699938e-12"
$ws.Range("I4").Value = [double]"2.754685368699938e-12"
$ws.Range("L4").Value = [double]"63.54186372250847"
$ws.Range("M4").Value = "[46.46218972920506, 80.62153771581188]"
$ws.Range("N4").Value = [double]"1.903163404648467e-09"
$ws.Range("O4").Value = [double]"1.903163404648467e-09"
$ws.Range("P4").Value = [double]"2.19502669890458"
$ws.Range("Q4").Value = "[1.9182898084982716, 2.4717635893108882]"
$ws.Range("R4").Value = [double]"0"
$ws.Range("S4").Value = [double]"0"
$ws.Range("T4").Value = [double]"67.81752837189451"
$ws.Range("U4").Value = "[58.8839228510611, 76.75113389272792]"
$ws.Range("V4").Value = [double]"0"
$ws.Range("W4").Value = [double]"0"
$ws.Range("X4").Value = [double]"15.13413413413426"
$ws.Range("Y4").Value = [double]"14.10966966966979"
$ws.Range("Z4").Value = [double]"16.15859859859874"
$ws.Range("F5").Value = [double]"23.2600000000002"
$ws.Range("H5").Value = [double]"1.867842902569805e-08"
$ws.Range("I5").Value = [double]"1.867842902569805e-08"
$ws.Range("L5").Value = [double]"57.28393308776188"
$ws.Range("M5").Value = "[35.85853350450074, 78.70933267102302]"
$ws.Range("N5").Value = [double]"2.528110724764332e-06"
$ws.Range("O5").Value = [double]"2.528110724764332e-06"
$ws.Range("P5").Value = [double]"2.207605648468503"
$ws.Range("Q5").Value = "[1.8176582119868865, 2.5975530849501194]"
$ws.Range("R5").Value = [double]"7.327471962526033e-15"
$ws.Range("S5").Value = [double]"7.327471962526033e-15"
$ws.Range("T5").Value = [double]"58.83023660705443"
$ws.Range("U5").Value = "[47.78657790613076, 69.8738953079781]"
$ws.Range("V5").Value = [double]"5.484501741648273e-14"
$ws.Range("W5").Value = [double]"5.484501741648273e-14"
$ws.Range("X5").Value = [double]"15.0875675675677"
$ws.Range("Y5").Value = [double]"13.64400400400412"
$ws.Range("Z5").Value = [double]"16.53113113113127"
$ws.Range("F6").Value = [double]"23.2600000000002"
$ws.Range("H6").Value = [double]"9.086311114625545e-08"
$ws.Range("I6").Value = [double]"9.086311114625545e-08"
$ws.Range("L6").Value = [double]"56.03564227250864"
$ws.Range("M6").Value = "[32.88827955583668, 79.1830049891806]"
$ws.Range("N6").Value = [double]"1.389512718197317e-05"
$ws.Range("O6").Value = [double]"1.389512718197317e-05"
$ws.Range("P6").Value = [double]"2.345974093671657"
$ws.Range("Q6").Value = "[1.9686056067539632, 2.72334258058935]"
$ws.Range("R6").Value = [double]"2.220446049250313e-16"
$ws.Range("S6").Value = [double]"2.220446049250313e-16"
$ws.Range("T6").Value = [double]"54.5841659357698"
$ws.Range("U6").Value = "[43.008598353263835, 66.15973351827576]"
$ws.Range("V6").Value = [double]"2.54996024295906e-12"
$ws.Range("W6").Value = [double]"2.54996024295906e-12"
$ws.Range("X6").Value = [double]"14.57533533533546"
$ws.Range("Y6").Value = [double]"13.17833833833845"
$ws.Range("Z6").Value = [double]"15.97233233233247"
$ws.Range("F7").Value = [double]"23.2600000000002"
$ws.Range("H7").Value = [double]"2.397522969044275e-07"
$ws.Range("I7").Value = [double]"2.397522969044275e-07"
$ws.Range("L7").Value = [double]"49.51235457950482"
$ws.Range("M7").Value = "[29.702907750539595, 69.32180140847004]"
$ws.Range("N7").Value = [double]"8.210644464767825e-06"
$ws.Range("O7").Value = [double]"8.210644464767825e-06"
$ws.Range("P7").Value = [double]"2.647868883205811"
$ws.Range("Q7").Value = "[2.220184598032426, 3.075553168379196]"
$ws.Range("R7").Value = [double]"2.220446049250313e-16"
$ws.Range("S7").Value = [double]"2.220446049250313e-16"
$ws.Range("T7").Value = [double]"52.38654579955399"
$ws.Range("U7").Value = "[41.60728087177495, 63.16581072733304]"
$ws.Range("V7").Value = [double]"1.011857264643368e-12"
$ws.Range("W7").Value = [double]"1.011857264643368e-12"
$ws.Range("X7").Value = [double]"13.45773773773785"
$ws.Range("Y7").Value = [double]"11.87447447447458"
$ws.Range("Z7").Value = [double]"15.04100100100113"
$ws.Range("F8").Value = [double]"23.2600000000002"
$ws.Range("H8").Value = [double]"2.456852354892902e-08"
$ws.Range("I8").Value = [double]"2.456852354892902e-08"
$ws.Range("L8").Value = [double]"49.80254184817222"
$ws.Range("M8").Value = "[32.22137752817697, 67.38370616816746]"
$ws.Range("N8").Value = [double]"8.53000112988056e-07"
$ws.Range("O8").Value = [double]"8.53000112988056e-07"
$ws.Range("P8").Value = [double]"2.849132076228581"
$ws.Range("Q8").Value = "[2.471763589310888, 3.2265005631462746]"
$ws.Range("T8").Value = [double]"56.2141922902656"
$ws.Range("U8").Value = "[46.363197505784264, 66.06518707474694]"
$ws.Range("V8").Value = [double]"5.551115123125783e-15"
$ws.Range("W8").Value = [double]"5.551115123125783e-15"
$ws.Range("X8").Value = [double]"12.71267267267278"
$ws.Range("Y8").Value = [double]"11.31567567567577"
$ws.Range("Z8").Value = [double]"14.10966966966979"
$ws.Range("F9").Value = [double]"22.81000000000013"
$ws.Range("H9").Value = [double]"2.402398280310081e-10"
$ws.Range("I9").Value = [double]"2.402398280310081e-10"
$ws.Range("L9").Value = [double]"63.08660208836336"
$ws.Range("M9").Value = "[46.82912405151919, 79.34408012520754]"
$ws.Range("N9").Value = [double]"6.40205666257998e-10"
$ws.Range("O9").Value = [double]"6.40205666257998e-10"
$ws.Range("P9").Value = [double]"-3.094421592725082"
$ws.Range("Q9").Value = "[-3.396316382259236, -2.792526803190927]"
$ws.Range("T9").Value = [double]"63.12762790852518"
$ws.Range("U9").Value = "[52.54753236852574, 73.70772344852463]"
$ws.Range("V9").Value = [double]"1.332267629550188e-15"
$ws.Range("W9").Value = [double]"1.332267629550188e-15"
$ws.Range("X9").Value = [double]"11.23375375375382"
$ws.Range("Y9").Value = [double]"10.13777777777783"
$ws.Range("Z9").Value = [double]"12.3297297297298"
$ws.Range("F10").Value = [double]"22.81000000000013"
$ws.Range("H10").Value = [double]"3.373022902097844e-07"
$ws.Range("I10").Value = [double]"3.373022902097844e-07"
$ws.Range("L10").Value = [double]"53.16782396340193"
$ws.Range("M10").Value = "[32.29630951228023, 74.03933841452363]"
$ws.Range("N10").Value = [double]"5.945867687318085e-06"
$ws.Range("O10").Value = [double]"5.945867687318085e-06"
$ws.Range("P10").Value = [double]"-2.855421551010543"
$ws.Range("Q10").Value = "[-3.3208426848756973, -2.390000417145388]"
$ws.Range("R10").Value = [double]"4.440892098500626e-16"
$ws.Range("S10").Value = [double]"4.440892098500626e-16"
$ws.Range("T10").Value = [double]"56.9121688304335"
$ws.Range("U10").Value = "[45.061002237134105, 68.7633354237329]"
$ws.Range("V10").Value = [double]"1.461719634221481e-12"
$ws.Range("W10").Value = [double]"1.461719634221481e-12"
$ws.Range("X10").Value = [double]"10.36610610610616"
$ws.Range("Y10").Value = [double]"8.676476476476523"
$ws.Range("Z10").Value = [double]"12.0557357357358"
$ws.Range("F11").Value = [double]"22.81000000000013"
$ws.Range("H11").Value = [double]"4.099019401682114e-08"
$ws.Range("I11").Value = [double]"4.099019401682114e-08"
$ws.Range("L11").Value = [double]"48.72261422472592"
$ws.Range("M11").Value = "[31.25893445638897, 66.18629399306286]"
$ws.Range("N11").Value = [double]"1.143389375179993e-06"
$ws.Range("O11").Value = [double]"1.143389375179993e-06"
$ws.Range("P11").Value = [double]"-2.754789954499158"
$ws.Range("Q11").Value = "[-3.1447373909807745, -2.3648425180175416]"
$ws.Range("T11").Value = [double]"53.4205540480076"
$ws.Range("U11").Value = "[43.55208882716755, 63.28901926884765]"
$ws.Range("V11").Value = [double]"3.241851231905457e-14"
$ws.Range("W11").Value = [double]"3.241851231905457e-14"
$ws.Range("X11").Value = [double]"10.00078078078084"
$ws.Range("Y11").Value = [double]"8.585145145145191"
$ws.Range("Z11").Value = [double]"11.41641641641648"
$ws.Range("F12").Value = [double]"22.81000000000013"
$ws.Range("H12").Value = [double]"0.0003024875562076534"
$ws.Range("I12").Value = [double]"0.0003024875562076534"
$ws.Range("L12").Value = [double]"36.54683915461952"
$ws.Range("M12").Value = "[14.737206872909404, 58.356471436329635]"
$ws.Range("N12").Value = [double]"0.001528329980261089"
$ws.Range("O12").Value = [double]"0.001528329980261089"
$ws.Range("P12").Value = [double]"-2.088105627611234"
$ws.Range("Q12").Value = "[-2.7044741562434655, -1.471737098979002]"
$ws.Range("R12").Value = [double]"1.863956966552394e-08"
$ws.Range("S12").Value = [double]"1.863956966552394e-08"
$ws.Range("T12").Value = [double]"51.12986429394139"
$ws.Range("U12").Value = "[39.53085600270565, 62.72887258517713]"
$ws.Range("V12").Value = [double]"1.883027067606236e-11"
$ws.Range("W12").Value = [double]"1.883027067606236e-11"
$ws.Range("X12").Value = [double]"7.580500500500543"
$ws.Range("Y12").Value = [double]"5.342882882882913"
$ws.Range("Z12").Value = [double]"9.818118118118173"
